# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.042.18"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "3.051.41"
$ws.Range("E3").Value = "  +5.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.16%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.433"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +5.08%  "
$ws.Range("E11").Value = "  +7.28%  "
$ws.Range("D12").Value = "3.562.73"
$ws.Range("E12").Value = "  +5.38%  "
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("D16").Value = "57.103.68"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "3.044.99"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "335.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.04%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.24%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "0.0₃0933"
$ws.Range("E27").Value = "  +12.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("E34").Value = "  +4.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.73%  "
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("D39").Value = "3.086.79"
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +6.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.663"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.21%  "
$ws.Range("D44").Value = "2.230.57"
$ws.Range("E44").Value = "  +7.19%  "
$ws.Range("E45").Value = "  +11.41%  "
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.940"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.45%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.47%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0869"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.684"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.52%  "
